$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 41: task "Magic Button para borrar movimientos de stock" is now complete (100%).
# Replace the "en proceso" text in C41 with a 100% numeric value (percentage style).
$ws.Range("C41").Value = 1
$ws.Range("C41").NumberFormat = "0%"

# Row 44: task "en stock agregar codigo - descripcion" gets assigned to Agustina
# and marked as "en proceso" (in progress).
$ws.Range("B44").Value = "Agustina"
$ws.Range("C44").Value = "en proceso"

# Move the active selection to C45 to reflect the new working position.
$ws.Range("C45").Select()
